$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.798.27"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "'1.859.70"
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("E4").Value = "  +0.47%  "
$ws.Range("D5").Value = "'323.50"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("D7").Value = "'0.4415"
$ws.Range("E7").Value = "  +0.75%  "
$ws.Range("D9").Value = "'0.07448"
$ws.Range("E9").Value = "  +0.35%  "
$ws.Range("D10").Value = "'0.8880"
$ws.Range("E10").Value = "  +1.28%  "
$ws.Range("D11").Value = "'21.64"
$ws.Range("E11").Value = "  +0.57%  "
$ws.Range("D12").Value = "'1.869.05"
$ws.Range("E12").Value = "  +0.60%  "
$ws.Range("D13").Value = "'5.548"
$ws.Range("E13").Value = "  +0.48%  "
$ws.Range("D14").Value = "'6.747"
$ws.Range("E14").Value = "  +0.59%  "
$ws.Range("D15").Value = "'0.07218"
$ws.Range("E15").Value = "  +0.31%  "
$ws.Range("D16").Value = "'86.22"
$ws.Range("E16").Value = "  +3.86%  "
$ws.Range("D17").Value = "'1.040"
$ws.Range("E17").Value = "  +0.38%  "
$ws.Range("D18").Value = "'0.000009119"
$ws.Range("E18").Value = "  +0.62%  "
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("E20").Value = "  +0.80%  "
$ws.Range("D21").Value = "'27.820.31"
$ws.Range("E21").Value = "  +0.54%  "
$ws.Range("D22").Value = "'5.300"
$ws.Range("E22").Value = "  +0.44%  "
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("D24").Value = "'2.073.74"
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("D25").Value = "'2.068"
$ws.Range("E25").Value = "  +6.26%  "
$ws.Range("D26").Value = "'159.38"
$ws.Range("E26").Value = "  +1.14%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("E28").Value = "  +3.73%  "
$ws.Range("D29").Value = "'5.388"
$ws.Range("E29").Value = "  +1.74%  "
$ws.Range("D30").Value = "'118.96"
$ws.Range("E30").Value = "  +2.13%  "
$ws.Range("D31").Value = "'0.09115"
$ws.Range("E31").Value = "  +0.29%  "
$ws.Range("E32").Value = "  +0.86%  "
$ws.Range("D33").Value = "'0.7753"
$ws.Range("E33").Value = "  +0.72%  "
$ws.Range("D34").Value = "'3.025"
$ws.Range("E34").Value = "  +4.90%  "
$ws.Range("D35").Value = "'4.614"
$ws.Range("E35").Value = "  +1.96%  "
$ws.Range("D36").Value = "'1.035"
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("D37").Value = "'1.157"
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").Value = "'0.01987"
$ws.Range("E38").Value = "  +0.22%  "
$ws.Range("D39").Value = "'0.05323"
$ws.Range("E39").Value = "  +0.51%  "
$ws.Range("D40").Value = "'2.867"
$ws.Range("E40").Value = "  +1.41%  "
$ws.Range("D41").Value = "'0.5218"
$ws.Range("E41").Value = "  +0.60%  "
$ws.Range("D42").Value = "'6.976"
$ws.Range("E42").Value = "  +3.52%  "
$ws.Range("E43").Value = "  +0.35%  "
$ws.Range("D44").Value = "'8.821"
$ws.Range("E44").Value = "  +2.70%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'10.81"
$ws.Range("E45").Value = "  +1.72%  "
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").Value = "'111.20"
$ws.Range("E46").Value = "  +2.00%  "
$ws.Range("E47").Value = "  +0.40%  "
$ws.Range("D48").Value = "'0.06587"
$ws.Range("E48").Value = "  +2.98%  "
$ws.Range("E49").Value = "  +0.10%  "
$ws.Range("D50").Value = "'0.4735"
$ws.Range("E50").Value = "  +1.50%  "
$ws.Range("D51").Value = "'1.888"
$ws.Range("E51").Value = "  -0.20%  "
